{"js": "// \"Update : Perbaiki judul PI\"\n// The thesis title inside the paragraph that reads:\n//   ... dengan judul \"Pembuatan Web E-Learning Programming Menggunakan PHP\n//   Native dan Framework Bootstrap\".\n// gets corrected from \"Web\" to \"Website\". This is a single in-place edit\n// (placing the cursor right after \"Web\" and typing \"site\"), which is also\n// exactly where Word's \"_GoBack\" bookmark (tracking the last edit point)\n// ends up afterwards - it moves from the end of the unrelated\n// \"Sistematika Tulisan Ilmiah\" heading to this new edit location.\n\n// 1) Find the unique occurrence of the title and, within it, the word \"Web\"\n//    that must become \"Website\".\nconst titleMatches = context.document.body.search(\"Pembuatan Web E-Learning\", {\n  matchCase: true,\n});\ntitleMatches.load(\"text\");\nawait context.sync();\n\nif (titleMatches.items.length > 0) {\n  const webToken = titleMatches.items[0].search(\"Web\", { matchCase: true });\n  webToken.load(\"text\");\n  await context.sync();\n\n  if (webToken.items.length > 0) {\n    // Insert \"site\" right after \"Web\" -> \"Website\".\n    webToken.items[0].insertText(\"site\", \"End\");\n    await context.sync();\n  }\n}\n\n// 2) Relocate the \"_GoBack\" bookmark: remove it from its old position...\nconst existingGoBack = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nawait context.sync();\nif (!existingGoBack.isNullObject) {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// ...and re-create it right before \"dan Framework Bootstrap\", which is where\n// the author's cursor was left after retyping the title.\nconst tailMatches = context.document.body.search(\"dan Framework Bootstrap\", {\n  matchCase: true,\n});\ntailMatches.load(\"text\");\nawait context.sync();\n\nif (tailMatches.items.length > 0) {\n  const cursorSpot = tailMatches.items[0].getRange(\"Start\");\n  cursorSpot.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# \"Update : Perbaiki judul PI\"\n# The thesis title inside the sentence:\n#   ... dengan judul \"Pembuatan Web E-Learning Programming Menggunakan PHP\n#   Native dan Framework Bootstrap\".\n# is corrected from \"Web\" to \"Website\". This is a single in-place edit\n# (cursor placed right after \"Web\", typing \"site\"). Word's \"_GoBack\"\n# bookmark - which always tracks the most recent edit location - follows\n# that edit: it moves from the end of the unrelated \"Sistematika Tulisan\n# Ilmiah\" heading to the new edit spot inside the title.\n\n$d = $word.ActiveDocument\n\n# 1) Locate the unique occurrence of the title and fix \"Web\" -> \"Website\".\n$titleRange = $d.Content\n$titleRange.Find.ClearFormatting()\n$titleRange.Find.Text = \"Pembuatan Web E-Learning\"\n$titleRange.Find.MatchCase = $true\n$titleRange.Find.MatchWholeWord = $false\n$titleRange.Find.Wrap = 1\n$titleRange.Find.Execute() | Out-Null\n\nif ($titleRange.Find.Found) {\n    $titleStart = $titleRange.Start\n    # \"Pembuatan \" is 10 characters, \"Web\" is the next 3.\n    $webRange = $d.Range($titleStart + 10, $titleStart + 13)\n    $webRange.Collapse(0)  # wdCollapseEnd -> insertion point right after \"Web\"\n    $webRange.InsertAfter(\"site\")\n}\n\n# 2) Move the \"_GoBack\" bookmark to sit right before \"dan Framework\n#    Bootstrap\", i.e. where the author's cursor ended up after retyping the\n#    title. Adding a bookmark with an already-existing name re-points it,\n#    so the stale one (after \"Sistematika Tulisan Ilmiah\") disappears.\n$tailRange = $d.Content\n$tailRange.Find.ClearFormatting()\n$tailRange.Find.Text = \"dan Framework Bootstrap\"\n$tailRange.Find.MatchCase = $true\n$tailRange.Find.MatchWholeWord = $false\n$tailRange.Find.Wrap = 1\n$tailRange.Find.Execute() | Out-Null\n\nif ($tailRange.Find.Found) {\n    $cursorSpot = $d.Range($tailRange.Start, $tailRange.Start)\n    $d.Bookmarks.Add(\"_GoBack\", $cursorSpot)\n}\n"}
